$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out the now-unused columns C through I for rows 3-10 entirely
$ws.Range("C3:I10").Clear()

# Row 2 (D2:H2) keep their cell/style but lose their date values
$ws.Range("D2:H2").ClearContents()

# Update the changed hour totals
$ws.Range("B7").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 18

# Update the selection to match the target state
$ws.Range("D7").Select()
